# ---------------------------------------------------------------------------
# Reproduces:
#   1) Slide 16's table (3rd shape, a graphicFrame) switching to the built-in
#      table style {8AF063CB-C343-43DB-9AA3-2C32800642FE}.
#   2) The presentation's theme colour scheme (ppt/theme/theme1.xml, the
#      theme actually used by the slide master / every slide) being swapped
#      from the "Integral" palette to the default "Office" palette - this is
#      what applying the built-in "Office Theme" design from the Design tab
#      changes on this deck (font scheme / format scheme were already
#      identical between the two themes, only the 12 colour-scheme slots
#      differ).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style -----------------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{8AF063CB-C343-43DB-9AA3-2C32800642FE}")

# --- 2) Theme colours -----------------------------------------------------
# MsoThemeColorSchemeIndex order: Dark1, Light1, Dark2, Light2,
# Accent1..Accent6, Hyperlink, FollowedHyperlink.
# Target values come from the standard Office theme:
#   000000, FFFFFF, 44546A, E7E6E6, 5B9BD5, ED7D31,
#   A5A5A5, FFC000, 4472C4, 70AD47, 0563C1, 954F72
# (PowerPoint's RGB long uses 0x00BBGGRR.)
$officeThemeColors = @(
    0,        # dk1      000000
    16777215, # lt1      FFFFFF
    6968388,  # dk2      44546A
    15132391, # lt2      E7E6E6
    13998939, # accent1  5B9BD5
    3243501,  # accent2  ED7D31
    10855845, # accent3  A5A5A5
    49407,    # accent4  FFC000
    12874308, # accent5  4472C4
    4697456,  # accent6  70AD47
    12673797, # hlink    0563C1
    7491477   # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i - 1]
}
